$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# K_Coexistence sheet (A1:C4): update raw data values, move the selection.
# ---------------------------------------------------------------------------
$wsK = $wb.Worksheets.Item("K_Coexistence")

$wsK.Range("B2").Value = 143.07
$wsK.Range("C2").Value = 8.36

$wsK.Range("A3").Value = 2000
$wsK.Range("B3").Value = 178.66
$wsK.Range("C3").Value = 11.02

$wsK.Range("A4").Value = 5000
$wsK.Range("B4").Value = 149.32
$wsK.Range("C4").Value = 13.44

# ---------------------------------------------------------------------------
# G_Coexistence sheet (A1:E11): only the remembered selection changes.
# ---------------------------------------------------------------------------
$wsG = $wb.Worksheets.Item("G_Coexistence")
$wsG.Range("B5").Select()

# ---------------------------------------------------------------------------
# Competition_Coexistence sheet (A1:E9): update the B/C inputs; the D/E
# columns are formulas (=B-C / =B+C) and recalculate automatically.
# ---------------------------------------------------------------------------
$wsC = $wb.Worksheets.Item("Competition_Coexistence")

$wsC.Range("B2").Value = 172.26
$wsC.Range("C2").Value = 9.9

$wsC.Range("B3").Value = 247.12
$wsC.Range("C3").Value = 13.01

$wsC.Range("B4").Value = 199.24
$wsC.Range("C4").Value = 11.54

$wsC.Range("B5").Value = 171.76
$wsC.Range("C5").Value = 11.19

$wsC.Range("B6").Value = 125.91
$wsC.Range("C6").Value = 8.53

$wsC.Range("B7").Value = 115.84
$wsC.Range("C7").Value = 7.69

$wsC.Range("B8").Value = 73.06
$wsC.Range("C8").Value = 4.13

$wsC.Range("B9").Value = 47.375
$wsC.Range("C9").Value = 2.41

$wsC.Range("F24").Select()

# ---------------------------------------------------------------------------
# Finally, activate K_Coexistence and select D15 - it becomes the workbook's
# active tab / active cell, matching the saved view state.
# ---------------------------------------------------------------------------
$wsK.Activate()
$wsK.Range("D15").Select()
